# Insert a new weekly record at row 15 (pushing existing rows 15:33 down to 16:34)
# and populate it with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("15:15").Insert()

$ws.Range("A15").Value = 10
$ws.Range("B15").Value = "Vega Modelo de Temuco"
$ws.Range("C15").Value = "La Araucanía"
$ws.Range("D15").Value = 44512
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = 100114002
$ws.Range("G15").Value = "Camote"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 30
$ws.Range("K15").Value = 20000
$ws.Range("L15").Value = 20000
$ws.Range("M15").Value = 20000
$ws.Range("N15").Value = '$/malla 20 kilos'
$ws.Range("O15").Value = "Perú"
$ws.Range("P15").Value = 1000
$ws.Range("Q15").Value = 20
$ws.Range("R15").Value = "Hortaliza"
